# Auto-generated edit script: update market/profit values on Anima Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# Sheet ALC row 62 (Leve Item ID 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2622.1875
$ws.Cells.Item(62, 9).Value = 1695
$ws.Cells.Item(62, 10).Value = 3814.2856
$ws.Cells.Item(62, 11).Value = 1695
$ws.Cells.Item(62, 12).Value = 3814.2856
$ws.Cells.Item(62, 13).Value = -1071
$ws.Cells.Item(62, 14).Value = -5062.2856

# Sheet ALC row 65 (Leve Item ID 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2622.1875
$ws.Cells.Item(65, 9).Value = 1695
$ws.Cells.Item(65, 10).Value = 3814.2856
$ws.Cells.Item(65, 11).Value = 8475
$ws.Cells.Item(65, 12).Value = 19071.428
$ws.Cells.Item(65, 13).Value = -5355
$ws.Cells.Item(65, 14).Value = -25311.428

# Sheet ALC row 96 (Leve Item ID 19894)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 1450.8
$ws.Cells.Item(96, 9).Value = 576.4
$ws.Cells.Item(96, 10).Value = 2325.2
$ws.Cells.Item(96, 11).Value = 1729.2
$ws.Cells.Item(96, 12).Value = 6975.599999999999
$ws.Cells.Item(96, 13).Value = -356.1999999999998
$ws.Cells.Item(96, 14).Value = -9721.599999999999

# Sheet ALC row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2351.087
$ws.Cells.Item(100, 9).Value = 2398.077
$ws.Cells.Item(100, 10).Value = 2290
$ws.Cells.Item(100, 11).Value = 2398.077
$ws.Cells.Item(100, 12).Value = 2290
$ws.Cells.Item(100, 13).Value = -1857.077
$ws.Cells.Item(100, 14).Value = -3372

# Sheet ALC row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2915.75
$ws.Cells.Item(113, 9).Value = 2897.8
$ws.Cells.Item(113, 10).Value = 2928.5715
$ws.Cells.Item(113, 11).Value = 2897.8
$ws.Cells.Item(113, 12).Value = 2928.5715
$ws.Cells.Item(113, 13).Value = 356.1999999999998
$ws.Cells.Item(113, 14).Value = -9436.5715

# Sheet ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 5643.478
$ws.Cells.Item(132, 9).Value = 5666.7144
$ws.Cells.Item(132, 10).Value = 5399.5
$ws.Cells.Item(132, 11).Value = 17000.1432
$ws.Cells.Item(132, 12).Value = 16198.5
$ws.Cells.Item(132, 13).Value = -14470.1432
$ws.Cells.Item(132, 14).Value = -21258.5

# Sheet ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2158.4036
$ws.Cells.Item(138, 10).Value = 2044.8379
$ws.Cells.Item(138, 12).Value = 6134.5137
$ws.Cells.Item(138, 14).Value = -16414.5137

# Sheet ARM row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3381.818
$ws.Cells.Item(61, 9).Value = 2360
$ws.Cells.Item(61, 10).Value = 4233.3335
$ws.Cells.Item(61, 11).Value = 2360
$ws.Cells.Item(61, 12).Value = 4233.3335
$ws.Cells.Item(61, 13).Value = -2148
$ws.Cells.Item(61, 14).Value = -4657.3335

# Sheet ARM row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 833.2
$ws.Cells.Item(74, 9).Value = 770.1818
$ws.Cells.Item(74, 11).Value = 770.1818
$ws.Cells.Item(74, 13).Value = 103.8182

# Sheet ARM row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 833.2
$ws.Cells.Item(77, 9).Value = 770.1818
$ws.Cells.Item(77, 11).Value = 3850.909
$ws.Cells.Item(77, 13).Value = 517.0910000000003

# Sheet ARM row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1745.3334
$ws.Cells.Item(110, 9).Value = 1494.4
$ws.Cells.Item(110, 11).Value = 1494.4
$ws.Cells.Item(110, 13).Value = 550.5999999999999

# Sheet ARM row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2105.6785
$ws.Cells.Item(122, 9).Value = 1769.591
$ws.Cells.Item(122, 10).Value = 3338
$ws.Cells.Item(122, 11).Value = 5308.772999999999
$ws.Cells.Item(122, 12).Value = 10014
$ws.Cells.Item(122, 13).Value = -2858.772999999999
$ws.Cells.Item(122, 14).Value = -14914

# Sheet ARM row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3381.818
$ws.Cells.Item(136, 9).Value = 2360
$ws.Cells.Item(136, 10).Value = 4233.3335
$ws.Cells.Item(136, 11).Value = 7080
$ws.Cells.Item(136, 12).Value = 12700.0005
$ws.Cells.Item(136, 13).Value = -4530
$ws.Cells.Item(136, 14).Value = -17800.0005

# Sheet CRP row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 883
$ws.Cells.Item(16, 9).Value = 877.75
$ws.Cells.Item(16, 10).Value = 890
$ws.Cells.Item(16, 11).Value = 877.75
$ws.Cells.Item(16, 12).Value = 890
$ws.Cells.Item(16, 13).Value = -590.75
$ws.Cells.Item(16, 14).Value = -1464

# Sheet CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2402.0278
$ws.Cells.Item(31, 9).Value = 969.1667
$ws.Cells.Item(31, 10).Value = 5267.75
$ws.Cells.Item(31, 11).Value = 969.1667
$ws.Cells.Item(31, 12).Value = 5267.75
$ws.Cells.Item(31, 13).Value = -674.1667
$ws.Cells.Item(31, 14).Value = -5857.75

# Sheet CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2402.0278
$ws.Cells.Item(34, 9).Value = 969.1667
$ws.Cells.Item(34, 10).Value = 5267.75
$ws.Cells.Item(34, 11).Value = 969.1667
$ws.Cells.Item(34, 12).Value = 5267.75
$ws.Cells.Item(34, 13).Value = -767.1667
$ws.Cells.Item(34, 14).Value = -5671.75

# Sheet CRP row 94 (Leve Item ID 32934)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1945
$ws.Cells.Item(94, 10).Value = 1945
$ws.Cells.Item(94, 12).Value = 1945
$ws.Cells.Item(94, 14).Value = -2847

# Sheet CRP row 98 (Leve Item ID 35357)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(98, 8).Value = 48799
$ws.Cells.Item(98, 10).Value = 48799
$ws.Cells.Item(98, 12).Value = 48799
$ws.Cells.Item(98, 14).Value = -53291

# Sheet CRP row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 883
$ws.Cells.Item(113, 9).Value = 877.75
$ws.Cells.Item(113, 10).Value = 890
$ws.Cells.Item(113, 11).Value = 877.75
$ws.Cells.Item(113, 12).Value = 890
$ws.Cells.Item(113, 13).Value = 1292.25
$ws.Cells.Item(113, 14).Value = -5230

# Sheet CUL row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2062.4707
$ws.Cells.Item(5, 9).Value = 873.25
$ws.Cells.Item(5, 11).Value = 2619.75
$ws.Cells.Item(5, 13).Value = -2507.75

# Sheet CUL row 80 (Leve Item ID 12890)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).ClearContents()

# Sheet CUL row 83 (Leve Item ID 12890)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).ClearContents()

# Sheet CUL row 122 (Leve Item ID 36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 4611.5386
$ws.Cells.Item(122, 9).Value = 508.2143
$ws.Cells.Item(122, 10).Value = 9398.75
$ws.Cells.Item(122, 11).Value = 4573.928699999999
$ws.Cells.Item(122, 12).Value = 84588.75
$ws.Cells.Item(122, 13).Value = -2123.928699999999
$ws.Cells.Item(122, 14).Value = -89488.75

# Sheet CUL row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 965.5823
$ws.Cells.Item(131, 9).Value = 312
$ws.Cells.Item(131, 10).Value = 1009.7432
$ws.Cells.Item(131, 11).Value = 936
$ws.Cells.Item(131, 12).Value = 3029.2296
$ws.Cells.Item(131, 13).Value = 4104
$ws.Cells.Item(131, 14).Value = -13109.2296

# Sheet CUL row 132 (Leve Item ID 43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 3581.855
$ws.Cells.Item(132, 10).Value = 4371.6
$ws.Cells.Item(132, 12).Value = 39344.4
$ws.Cells.Item(132, 14).Value = -44404.4

# Sheet CUL row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 2062.4707
$ws.Cells.Item(135, 9).Value = 873.25
$ws.Cells.Item(135, 11).Value = 7859.25
$ws.Cells.Item(135, 13).Value = -5324.25

# Sheet CUL row 137 (Leve Item ID 44088)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 8383.0625
$ws.Cells.Item(137, 10).Value = 3000
$ws.Cells.Item(137, 12).Value = 9000
$ws.Cells.Item(137, 14).Value = -19200

# Sheet CUL row 140 (Leve Item ID 44097)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 1984.96
$ws.Cells.Item(140, 9).Value = 1243.1765
$ws.Cells.Item(140, 10).Value = 3561.25
$ws.Cells.Item(140, 11).Value = 3729.5295
$ws.Cells.Item(140, 12).Value = 10683.75
$ws.Cells.Item(140, 13).Value = 1450.4705
$ws.Cells.Item(140, 14).Value = -21043.75

# Sheet GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4089.487
$ws.Cells.Item(122, 9).Value = 1898.8462
$ws.Cells.Item(122, 10).Value = 5184.8076
$ws.Cells.Item(122, 11).Value = 5696.5386
$ws.Cells.Item(122, 12).Value = 15554.4228
$ws.Cells.Item(122, 13).Value = -3246.5386
$ws.Cells.Item(122, 14).Value = -20454.4228

# Sheet GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3708.25
$ws.Cells.Item(132, 9).Value = 3900
$ws.Cells.Item(132, 10).Value = 3571.2856
$ws.Cells.Item(132, 11).Value = 11700
$ws.Cells.Item(132, 12).Value = 10713.8568
$ws.Cells.Item(132, 13).Value = -9170
$ws.Cells.Item(132, 14).Value = -15773.8568

# Sheet LTW row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 416.96295
$ws.Cells.Item(55, 9).Value = 260.5
$ws.Cells.Item(55, 10).Value = 509
$ws.Cells.Item(55, 11).Value = 260.5
$ws.Cells.Item(55, 12).Value = 509
$ws.Cells.Item(55, 13).Value = -87.5
$ws.Cells.Item(55, 14).Value = -855

# Sheet LTW row 61 (Leve Item ID 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4970
$ws.Cells.Item(61, 9).Value = 5458.8887
$ws.Cells.Item(61, 10).Value = 4090
$ws.Cells.Item(61, 11).Value = 5458.8887
$ws.Cells.Item(61, 12).Value = 4090
$ws.Cells.Item(61, 13).Value = -5256.8887
$ws.Cells.Item(61, 14).Value = -4494

# Sheet LTW row 113 (Leve Item ID 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 4970
$ws.Cells.Item(113, 9).Value = 5458.8887
$ws.Cells.Item(113, 10).Value = 4090
$ws.Cells.Item(113, 11).Value = 5458.8887
$ws.Cells.Item(113, 12).Value = 4090
$ws.Cells.Item(113, 13).Value = -3288.8887
$ws.Cells.Item(113, 14).Value = -8430

# Sheet WVR row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1294.2646
$ws.Cells.Item(126, 9).Value = 974.0476
$ws.Cells.Item(126, 10).Value = 1811.5385
$ws.Cells.Item(126, 11).Value = 2922.1428
$ws.Cells.Item(126, 12).Value = 5434.6155
$ws.Cells.Item(126, 13).Value = -452.1428000000001
$ws.Cells.Item(126, 14).Value = -10374.6155

